# Update NATMI TPM-derived statistics (ligand/receptor specificity + edge
# weights) for the Fn1-Mag LR-pair sheet to reflect the refreshed TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> MuSCs)
$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 0.232947
$ws.Range("N2").Value = 0.698841
$ws.Range("O2").Value = 0.2572219815457369
$ws.Range("P2").Value = 0.2572219815457369
$ws.Range("Q2").Value = 6.807514075362
$ws.Range("R2").Value = 61.267626678258
$ws.Range("S2").Value = 0.004348112379406882
$ws.Range("T2").Value = 0.004348112379406882

# Row 3 (ECs -> Resolving-Mac)
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("M3").Value = 0.6726793333333333
$ws.Range("N3").Value = 2.018038
$ws.Range("O3").Value = 0.7427780184542632
$ws.Range("P3").Value = 0.7427780184542632
$ws.Range("Q3").Value = 19.65800817298267
$ws.Range("R3").Value = 176.922073556844
$ws.Range("S3").Value = 0.01255601203981092
$ws.Range("T3").Value = 0.01255601203981092

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.9471112884046843
$ws.Range("J4").Value = 0.9471112884046842
$ws.Range("M4").Value = 0.232947
$ws.Range("N4").Value = 0.698841
$ws.Range("O4").Value = 0.2572219815457369
$ws.Range("P4").Value = 0.2572219815457369
$ws.Range("Q4").Value = 381.41421979947
$ws.Range("R4").Value = 3432.72797819523
$ws.Range("S4").Value = 0.2436178423477888
$ws.Range("T4").Value = 0.2436178423477888

# Row 5 (FAPs -> Resolving-Mac)
$ws.Range("I5").Value = 0.9471112884046843
$ws.Range("J5").Value = 0.9471112884046842
$ws.Range("M5").Value = 0.6726793333333333
$ws.Range("N5").Value = 2.018038
$ws.Range("O5").Value = 0.7427780184542632
$ws.Range("P5").Value = 0.7427780184542632
$ws.Range("Q5").Value = 1101.407028631238
$ws.Range("R5").Value = 9912.663257681139
$ws.Range("S5").Value = 0.7034934460568956
$ws.Range("T5").Value = 0.7034934460568955

# Row 6 (MuSCs -> MuSCs)
$ws.Range("G6").Value = 37.39212666666667
$ws.Range("H6").Value = 112.17638
$ws.Range("I6").Value = 0.02162924801792661
$ws.Range("J6").Value = 0.0216292480179266
$ws.Range("M6").Value = 0.232947
$ws.Range("N6").Value = 0.698841
$ws.Range("O6").Value = 0.2572219815457369
$ws.Range("P6").Value = 0.2572219815457369
$ws.Range("Q6").Value = 8.710383730620002
$ws.Range("R6").Value = 78.39345357558001
$ws.Range("S6").Value = 0.005563518034515284
$ws.Range("T6").Value = 0.005563518034515283

# Row 7 (MuSCs -> Resolving-Mac)
$ws.Range("G7").Value = 37.39212666666667
$ws.Range("H7").Value = 112.17638
$ws.Range("I7").Value = 0.02162924801792661
$ws.Range("J7").Value = 0.0216292480179266
$ws.Range("M7").Value = 0.6726793333333333
$ws.Range("N7").Value = 2.018038
$ws.Range("O7").Value = 0.7427780184542632
$ws.Range("P7").Value = 0.7427780184542632
$ws.Range("Q7").Value = 25.15291083804889
$ws.Range("R7").Value = 226.37619754244
$ws.Range("S7").Value = 0.01606572998341133
$ws.Range("T7").Value = 0.01606572998341132

# Row 8 (Resolving-Mac -> MuSCs)
$ws.Range("G8").Value = 24.817167
$ws.Range("H8").Value = 74.45150100000001
$ws.Range("I8").Value = 0.01435533915817136
$ws.Range("J8").Value = 0.01435533915817136
$ws.Range("M8").Value = 0.232947
$ws.Range("N8").Value = 0.698841
$ws.Range("O8").Value = 0.2572219815457369
$ws.Range("P8").Value = 0.2572219815457369
$ws.Range("Q8").Value = 5.781084601149001
$ws.Range("R8").Value = 52.02976141034101
$ws.Range("S8").Value = 0.003692508784025948
$ws.Range("T8").Value = 0.003692508784025948

# Row 9 (Resolving-Mac -> Resolving-Mac)
$ws.Range("G9").Value = 24.817167
$ws.Range("H9").Value = 74.45150100000001
$ws.Range("I9").Value = 0.01435533915817136
$ws.Range("J9").Value = 0.01435533915817136
$ws.Range("M9").Value = 0.6726793333333333
$ws.Range("N9").Value = 2.018038
$ws.Range("O9").Value = 0.7427780184542632
$ws.Range("P9").Value = 0.7427780184542632
$ws.Range("Q9").Value = 16.693995352782
$ws.Range("R9").Value = 150.245958175038
$ws.Range("S9").Value = 0.01066283037414541
$ws.Range("T9").Value = 0.01066283037414541
